$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organograma")
Write-Output ("before=" + $ws.AutoFilter.Range.Address())
$ws.AutoFilterMode = $false
$ws.Range("A1:E35").AutoFilter() | Out-Null
Write-Output ("after=" + $ws.AutoFilter.Range.Address())
